$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "MuSCs" rows (8 and 9) are removed from the data entirely, and the
# "Inflammatory-Mac" cluster rows (6 and 7) are re-labelled as "MuSCs" using
# updated TPM-derived values. Delete the higher-numbered row first so the
# second delete still targets the intended row.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

$ws.Cells.Item(6, 1).Value2 = "MuSCs"
$ws.Cells.Item(7, 1).Value2 = "MuSCs"

# Row 2 (ECs -> ECs)
$ws.Cells.Item(2, 9).Value2  = 0.01683020325561672     # I2
$ws.Cells.Item(2, 10).Value2 = 0.01683020325561672     # J2
$ws.Cells.Item(2, 13).Value2 = 0.08241233333333334     # M2
$ws.Cells.Item(2, 17).Value2 = 0.004404664508888889    # Q2
$ws.Cells.Item(2, 18).Value2 = 0.03964198058000001     # R2
$ws.Cells.Item(2, 19).Value2 = 0.0007955646525490005   # S2
$ws.Cells.Item(2, 20).Value2 = 0.0007955646525490005   # T2

# Row 3 (ECs -> FAPs)
$ws.Cells.Item(3, 9).Value2  = 0.01683020325561672     # I3
$ws.Cells.Item(3, 10).Value2 = 0.01683020325561672     # J3
$ws.Cells.Item(3, 19).Value2 = 0.01603463860306772     # S3
$ws.Cells.Item(3, 20).Value2 = 0.01603463860306771     # T3

# Row 4 (FAPs -> ECs)
$ws.Cells.Item(4, 9).Value2  = 0.9501703593606328      # I4
$ws.Cells.Item(4, 10).Value2 = 0.9501703593606329      # J4
$ws.Cells.Item(4, 13).Value2 = 0.08241233333333334     # M4
$ws.Cells.Item(4, 17).Value2 = 0.2486708921876667      # Q4
$ws.Cells.Item(4, 19).Value2 = 0.04491460621872336     # S4
$ws.Cells.Item(4, 20).Value2 = 0.04491460621872337     # T4

# Row 5 (FAPs -> FAPs)
$ws.Cells.Item(5, 9).Value2  = 0.9501703593606328      # I5
$ws.Cells.Item(5, 10).Value2 = 0.9501703593606329      # J5
$ws.Cells.Item(5, 19).Value2 = 0.9052557531419094      # S5
$ws.Cells.Item(5, 20).Value2 = 0.9052557531419094      # T5

# Row 6 (now MuSCs -> ECs)
$ws.Cells.Item(6, 7).Value2  = 0.1047943333333333      # G6
$ws.Cells.Item(6, 8).Value2  = 0.314383                # H6
$ws.Cells.Item(6, 9).Value2  = 0.03299943738375047     # I6
$ws.Cells.Item(6, 10).Value2 = 0.03299943738375048     # J6
$ws.Cells.Item(6, 13).Value2 = 0.08241233333333334     # M6
$ws.Cells.Item(6, 17).Value2 = 0.008636345530111111    # Q6
$ws.Cells.Item(6, 18).Value2 = 0.07772710977100002     # R6
$ws.Cells.Item(6, 19).Value2 = 0.001559885257342599    # S6
$ws.Cells.Item(6, 20).Value2 = 0.0015598852573426      # T6

# Row 7 (now MuSCs -> FAPs)
$ws.Cells.Item(7, 7).Value2  = 0.1047943333333333      # G7
$ws.Cells.Item(7, 8).Value2  = 0.314383                # H7
$ws.Cells.Item(7, 9).Value2  = 0.03299943738375047     # I7
$ws.Cells.Item(7, 10).Value2 = 0.03299943738375048     # J7
$ws.Cells.Item(7, 17).Value2 = 0.1740659027306667      # Q7
$ws.Cells.Item(7, 18).Value2 = 1.566593124576          # R7
$ws.Cells.Item(7, 19).Value2 = 0.03143955212640787     # S7
$ws.Cells.Item(7, 20).Value2 = 0.03143955212640787     # T7
